$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.303.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.82%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.095.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.16%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "387.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.35"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.50%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.86%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.68%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.89"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.14%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0855"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.577.20"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.90%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.20%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.73"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.095.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.07%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.990"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.24%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.64"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.404.17"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.59%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.23"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.49"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.05%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.18%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.94"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.28%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.60"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.93%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.15"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.08%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.96"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.42%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.25"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.93%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.80%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.16%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.166"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.51%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.35"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.98"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.33%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0472"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.20%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.46%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.81"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.36%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.78%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.290"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "130.60"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.39%  "

# Row 42
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.23%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.115"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.51%  "

# Row 44
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.52"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.05%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.22%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.59%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.17%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.071.45"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.08%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +18.22%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.28%  "
